# Update "想去人数" (number of people who want to go) counts on the
# "展览" and "全部类型" sheets to reflect the latest scrape snapshot.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 7508
$wsExhibit.Range("F5").Value = 9
$wsExhibit.Range("F6").Value = 451
$wsExhibit.Range("F7").Value = 4109
$wsExhibit.Range("F10").Value = 275
$wsExhibit.Range("F11").Value = 657
$wsExhibit.Range("F12").Value = 145

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 7508
$wsAll.Range("F7").Value = 9
$wsAll.Range("F8").Value = 451
$wsAll.Range("F9").Value = 4109
$wsAll.Range("F12").Value = 275
$wsAll.Range("F13").Value = 657
$wsAll.Range("F15").Value = 145
